$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column O into the new column P for every row that has data in O
# (row 15 has no O value, so it is skipped to avoid creating a stray empty P15 cell)
$ws.Range("O4:O14").Copy() | Out-Null
$ws.Range("P4:P14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("O16:O17").Copy() | Out-Null
$ws.Range("P16:P17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new 2022 column (P) values
$ws.Range("P4").Value = 2022

$ws.Range("P5").Value = 11.4
$ws.Range("P6").Value = 12.6
$ws.Range("P7").Value = 9.8000000000000007
$ws.Range("P8").Value = 11.4
$ws.Range("P9").Value = 5.4
$ws.Range("P10").Value = 4.7
$ws.Range("P11").Value = 3.4
$ws.Range("P12").Value = 17.7
$ws.Range("P13").Value = 20.5
$ws.Range("P14").Value = 8.4

$ws.Range("P16").Value = 12.9
$ws.Range("P17").Value = 10.5

# Update the saved selection to match the new active cell
$ws.Range("Q4").Select() | Out-Null
